$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.03448393990524
$ws.Range("D2").Value = 1.043749958186429
$ws.Range("E2").Value = 1.051573381290483
$ws.Range("F2").Value = 1.056725289445023
$ws.Range("I2").Value = 1.036513421144449
$ws.Range("J2").Value = 1.03960293354099
$ws.Range("K2").Value = 1.046522945324339
$ws.Range("L2").Value = 1.05432450838179
$ws.Range("M2").Value = 1.059462212687779
$ws.Range("N2").Value = 1.041079288950741
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.035490721134951
$ws.Range("D3").Value = 1.044554239357529
$ws.Range("E3").Value = 1.052565616933042
$ws.Range("F3").Value = 1.05772117011777
$ws.Range("I3").Value = 1.03671883634696
$ws.Range("J3").Value = 1.040252575667629
$ws.Range("K3").Value = 1.04713834262913
$ws.Range("L3").Value = 1.055128955602534
$ws.Range("M3").Value = 1.060271325536641
$ws.Range("N3").Value = 1.041729853643715
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.036142181717929
$ws.Range("D4").Value = 1.045074307479887
$ws.Range("E4").Value = 1.05320849848678
$ws.Range("F4").Value = 1.058365977726534
$ws.Range("I4").Value = 1.036849849364788
$ws.Range("J4").Value = 1.04067235798581
$ws.Range("K4").Value = 1.047535538723946
$ws.Range("L4").Value = 1.055649709693872
$ws.Range("M4").Value = 1.060794666604938
$ws.Range("N4").Value = 1.042150232100938
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.03641605683796
$ws.Range("D5").Value = 1.045292858037538
$ws.Range("E5").Value = 1.053478965849464
$ws.Range("F5").Value = 1.058637151250498
$ws.Range("I5").Value = 1.036904470803422
$ws.Range("J5").Value = 1.040848695157916
$ws.Range("K5").Value = 1.047702278119581
$ws.Range("L5").Value = 1.055868687368763
$ws.Range("M5").Value = 1.061014628382317
$ws.Range("N5").Value = 1.042326819692054
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.036462041730982
$ws.Range("D6").Value = 1.045329548562094
$ws.Range("E6").Value = 1.053524390210936
$ws.Range("F6").Value = 1.058682688107012
$ws.Range("I6").Value = 1.036913615202351
$ws.Range("J6").Value = 1.04087829475914
$ws.Range("K6").Value = 1.047730260194451
$ws.Range("L6").Value = 1.055905457750681
$ws.Range("M6").Value = 1.061051557931024
$ws.Range("N6").Value = 1.042356461328107
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.036145841246779
$ws.Range("D7").Value = 1.045077228102591
$ws.Range("E7").Value = 1.05321211170028
$ws.Range("F7").Value = 1.058369600782696
$ws.Range("I7").Value = 1.036850581013272
$ws.Range("J7").Value = 1.040674714759862
$ws.Range("K7").Value = 1.04753776765513
$ws.Range("L7").Value = 1.055652635478294
$ws.Range("M7").Value = 1.060797605945579
$ws.Range("N7").Value = 1.042152592221879
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.034824184954648
$ws.Range("D8").Value = 1.04402184182625
$ws.Range("E8").Value = 1.051908537800145
$ws.Range("F8").Value = 1.057061767363892
$ws.Range("I8").Value = 1.036583236167323
$ws.Range("J8").Value = 1.039822602808264
$ws.Range("K8").Value = 1.046731129812378
$ws.Range("L8").Value = 1.05459632849223
$ws.Range("M8").Value = 1.05973569899223
$ws.Range("N8").Value = 1.041299270173572
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.032495314752881
$ws.Range("D9").Value = 1.042159434639024
$ws.Range("E9").Value = 1.049617934489681
$ws.Range("F9").Value = 1.054760345725375
$ws.Range("I9").Value = 1.036097578206454
$ws.Range("J9").Value = 1.038316658122545
$ws.Range("K9").Value = 1.045302049603914
$ws.Range("L9").Value = 1.052736719503336
$ws.Range("M9").Value = 1.057862915292097
$ws.Range("N9").Value = 1.039791186873665
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.030942786911454
$ws.Range("D10").Value = 1.04091608069364
$ws.Range("E10").Value = 1.048095256712117
$ws.Range("F10").Value = 1.053228228388559
$ws.Range("I10").Value = 1.035764045073667
$ws.Range("J10").Value = 1.037309757696787
$ws.Range("K10").Value = 1.044344206156392
$ws.Range("L10").Value = 1.051498189770117
$ws.Range("M10").Value = 1.056613381115418
$ws.Range("N10").Value = 1.038782856533809
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.030270539580505
$ws.Range("D11").Value = 1.0403772905146
$ws.Range("E11").Value = 1.047436971306114
$ws.Range("F11").Value = 1.05256532710652
$ws.Range("I11").Value = 1.035617312171189
$ws.Range("J11").Value = 1.036873067860233
$ws.Range("K11").Value = 1.04392824221962
$ws.Range("L11").Value = 1.050962187572086
$ws.Range("M11").Value = 1.056072087346587
$ws.Range("N11").Value = 1.038345546547597
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.030020838154397
$ws.Range("D12").Value = 1.040177099078529
$ws.Range("E12").Value = 1.047192612351949
$ws.Range("F12").Value = 1.052319174184522
$ws.Range("I12").Value = 1.035562462264826
$ws.Range("J12").Value = 1.036710757741355
$ws.Range("K12").Value = 1.043773553257697
$ws.Range("L12").Value = 1.050763136535992
$ws.Range("M12").Value = 1.055870991687246
$ws.Range("N12").Value = 1.038183005929734
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.030074399904579
$ws.Range("D13").Value = 1.040220043558813
$ws.Range("E13").Value = 1.047245021036354
$ws.Range("F13").Value = 1.05237197127576
$ws.Range("I13").Value = 1.035574243441666
$ws.Range("J13").Value = 1.036745578515962
$ws.Range("K13").Value = 1.043806742766449
$ws.Range("L13").Value = 1.050805831656573
$ws.Range("M13").Value = 1.055914128965616
$ws.Range("N13").Value = 1.038217876153834
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.030249899158735
$ws.Range("D14").Value = 1.040360743863321
$ws.Range("E14").Value = 1.047416769300669
$ws.Range("F14").Value = 1.052544978431563
$ws.Range("I14").Value = 1.035612785332548
$ws.Range("J14").Value = 1.036859653375914
$ws.Range("K14").Value = 1.043915459279747
$ws.Range("L14").Value = 1.050945733049515
$ws.Range("M14").Value = 1.056055465428328
$ws.Range("N14").Value = 1.038332113013173
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.030358030196961
$ws.Range("D15").Value = 1.040447425934734
$ws.Range("E15").Value = 1.047522609972951
$ws.Range("F15").Value = 1.05265158422289
$ws.Range("I15").Value = 1.03563648632485
$ws.Range("J15").Value = 1.036929924883459
$ws.Range("K15").Value = 1.043982419086576
$ws.Range("L15").Value = 1.051031936768662
$ws.Range("M15").Value = 1.056142542860417
$ws.Range("N15").Value = 1.038402484314318
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.030987401913833
$ws.Range("D16").Value = 1.040951829857967
$ws.Range("E16").Value = 1.048138967003739
$ws.Range("F16").Value = 1.053272233883833
$ws.Range("I16").Value = 1.035773734600918
$ws.Range("J16").Value = 1.037338724747271
$ws.Range("K16").Value = 1.044371786838392
$ws.Range("L16").Value = 1.051533768616058
$ws.Range("M16").Value = 1.056649300052369
$ws.Range("N16").Value = 1.038811864720826
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.031382192194464
$ws.Range("D17").Value = 1.041268119974213
$ws.Range("E17").Value = 1.048525871512761
$ws.Range("F17").Value = 1.053661689537422
$ws.Range("I17").Value = 1.035859208425668
$ws.Range("J17").Value = 1.037594968226187
$ws.Range("K17").Value = 1.044615703015778
$ws.Range("L17").Value = 1.05184863234914
$ws.Range("M17").Value = 1.056967112355232
$ws.Range("N17").Value = 1.039068472094874
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.031612467403833
$ws.Range("D18").Value = 1.041452567064206
$ws.Range("E18").Value = 1.048751647063367
$ws.Range("F18").Value = 1.053888902138305
$ws.Range("I18").Value = 1.035908840870536
$ws.Range("J18").Value = 1.037744363598978
$ws.Range("K18").Value = 1.044757858274952
$ws.Range("L18").Value = 1.052032314917251
$ws.Range("M18").Value = 1.05715246396662
$ws.Range("N18").Value = 1.039218079626231
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.031690985419851
$ws.Range("D19").Value = 1.041515452017409
$ws.Range("E19").Value = 1.048828647750426
$ws.Range("F19").Value = 1.053966384174522
$ws.Range("I19").Value = 1.03592572641603
$ws.Range("J19").Value = 1.037795292164826
$ws.Range("K19").Value = 1.044806309695719
$ws.Range("L19").Value = 1.052094950637558
$ws.Range("M19").Value = 1.057215660169266
$ws.Range("N19").Value = 1.039269080516484
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.031339834858908
$ws.Range("D20").Value = 1.041234189126658
$ws.Range("E20").Value = 1.048484349900013
$ws.Range("F20").Value = 1.053619899485235
$ws.Range("I20").Value = 1.035850060955199
$ws.Range("J20").Value = 1.03756748266487
$ws.Range("K20").Value = 1.044589545207794
$ws.Range("L20").Value = 1.051814847562469
$ws.Range("M20").Value = 1.056933016470082
$ws.Range("N20").Value = 1.039040947500908
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.030198218970116
$ws.Range("D21").Value = 1.040319312801925
$ws.Range("E21").Value = 1.047366189362712
$ws.Range("F21").Value = 1.052494029974211
$ws.Range("I21").Value = 1.035601445275707
$ws.Range("J21").Value = 1.036826064042611
$ws.Range("K21").Value = 1.043883449980258
$ws.Range("L21").Value = 1.050904534355767
$ws.Range("M21").Value = 1.056013846325698
$ws.Range("N21").Value = 1.038298475979163
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.02948044650697
$ws.Range("D22").Value = 1.039743742507607
$ws.Range("E22").Value = 1.046664069361643
$ws.Range("F22").Value = 1.051786603486722
$ws.Range("I22").Value = 1.035443124606292
$ws.Range("J22").Value = 1.036359302802515
$ws.Range("K22").Value = 1.043438450434107
$ws.Range("L22").Value = 1.050332439430137
$ws.Range("M22").Value = 1.055435725059332
$ws.Range("N22").Value = 1.037831051884568
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.029860950616338
$ws.Range("D23").Value = 1.040048896245083
$ws.Range("E23").Value = 1.047036189791107
$ws.Range("F23").Value = 1.05216158055444
$ws.Range("I23").Value = 1.035527243411501
$ws.Range("J23").Value = 1.036606798693043
$ws.Range("K23").Value = 1.043674452402163
$ws.Range("L23").Value = 1.050635693401343
$ws.Range("M23").Value = 1.055742217081957
$ws.Range("N23").Value = 1.038078899247649
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.031358974313233
$ws.Range("D24").Value = 1.041249521138433
$ws.Range("E24").Value = 1.048503111418392
$ws.Range("F24").Value = 1.05363878245737
$ws.Range("I24").Value = 1.035854194992347
$ws.Range("J24").Value = 1.037579902414244
$ws.Range("K24").Value = 1.044601365156505
$ws.Range("L24").Value = 1.051830113367153
$ws.Range("M24").Value = 1.056948423003752
$ws.Range("N24").Value = 1.03905338488775
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.033097375314754
$ws.Range("D25").Value = 1.042641223754094
$ws.Range("E25").Value = 1.050209339852901
$ws.Range("F25").Value = 1.055354940448401
$ws.Range("I25").Value = 1.036224855050518
$ws.Range("J25").Value = 1.038706500876611
$ws.Range("K25").Value = 1.045672407096867
$ws.Range("L25").Value = 1.053217262753182
$ws.Range("M25").Value = 1.058347256182799
$ws.Range("N25").Value = 1.040181583249157
